$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 18) following the same layout as the existing rows
$row = 18

$ws.Range("A" + ($row - 1)).Copy($ws.Range("A" + $row))
$ws.Cells.Item($row, 1).Value = 42625.885057870371

$ws.Cells.Item($row, 2).Value = -26
$ws.Cells.Item($row, 3).Value = 63
$ws.Cells.Item($row, 4).Value = 35
$ws.Cells.Item($row, 5).Value = 80
$ws.Cells.Item($row, 6).Value = 20
$ws.Cells.Item($row, 7).Value = 6649
$ws.Cells.Item($row, 8).Value = 11170
$ws.Cells.Item($row, 9).Value = 1190
$ws.Cells.Item($row, 10).Value = 227
$ws.Cells.Item($row, 11).Value = 127
$ws.Cells.Item($row, 12).Value = 16
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = "Bag"
